$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 4 '61.331.36'
Set-TextValue 2 5 '  -0.31%  '
Set-TextValue 3 4 '3.381.84'
Set-TextValue 3 5 '  +1.90%  '
Set-TextValue 4 4 '1.00'
Set-TextValue 4 5 '  -0.04%  '
Set-TextValue 5 4 '572.60'
Set-TextValue 5 5 '  +0.61%  '
Set-TextValue 6 4 '137.10'
Set-TextValue 6 5 '  +8.08%  '
Set-TextValue 7 5 '  -0.04%  '
Set-TextValue 8 4 '3.382.17'
Set-TextValue 8 5 '  +1.96%  '
Set-TextValue 9 4 '0.475'
Set-TextValue 9 5 '  -0.14%  '
Set-TextValue 10 4 '7.60'
Set-TextValue 10 5 '  +5.03%  '
Set-TextValue 11 5 '  +4.74%  '
Set-TextValue 12 4 '0.391'
Set-TextValue 12 5 '  +4.26%  '
Set-TextValue 13 4 '3.956.24'
Set-TextValue 13 5 '  +1.81%  '
Set-TextValue 14 5 '  +2.38%  '
Set-TextValue 15 5 '  +2.17%  '
Set-TextValue 16 4 '3.391.22'
Set-TextValue 16 5 '  +2.35%  '
Set-TextValue 17 4 '25.28'
Set-TextValue 17 5 '  +2.48%  '
Set-TextValue 18 4 '61.423.51'
Set-TextValue 18 5 '  -0.24%  '
Set-TextValue 19 4 '14.02'
Set-TextValue 19 5 '  +6.72%  '
Set-TextValue 20 5 '  +4.41%  '
Set-TextValue 21 4 '9.43'
Set-TextValue 21 5 '  +3.32%  '
Set-TextValue 22 4 '379.66'
Set-TextValue 22 5 '  +6.17%  '
Set-TextValue 23 4 '0.574'
Set-TextValue 23 5 '  +4.02%  '
Set-TextValue 24 4 '3.513.12'
Set-TextValue 24 5 '  +1.81%  '
Set-TextValue 25 5 '  +0.09%  '
Set-TextValue 26 4 '70.98'
Set-TextValue 26 5 '  +0.75%  '
Set-TextValue 27 5 '  +10.24%  '
Set-TextValue 28 5 '  +12.03%  '
Set-TextValue 29 4 '7.76'
Set-TextValue 29 5 '  +7.65%  '
Set-TextValue 30 4 '0.999'
Set-TextValue 30 5 '  -0.07%  '
Set-TextValue 31 5 '  +3.24%  '
Set-TextValue 32 5 '  +4.89%  '
Set-TextValue 33 5 '  +2.05%  '
Set-TextValue 35 4 '3.410.58'
Set-TextValue 35 5 '  +1.85%  '
Set-TextValue 36 4 '23.55'
Set-TextValue 36 5 '  +5.37%  '
Set-TextValue 37 4 '5.56'
Set-TextValue 37 5 '  +0.59%  '
Set-TextValue 38 4 '6.98'
Set-TextValue 38 5 '  +4.24%  '
Set-TextValue 39 5 '  +3.95%  '
Set-TextValue 40 4 '163.54'
Set-TextValue 40 5 '  +0.04%  '
Set-TextValue 41 4 '0.0801'
Set-TextValue 41 5 '  +5.80%  '
Set-TextValue 42 4 '0.999'
Set-TextValue 42 5 '  -0.11%  '
Set-TextValue 43 2 'Filecoin'
Set-TextValue 43 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 43 4 '4.43'
Set-TextValue 43 5 '  +4.91%  '
Set-TextValue 44 2 'OKB'
Set-TextValue 44 3 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 44 4 '41.52'
Set-TextValue 44 5 '  +1.35%  '
Set-TextValue 45 2 'ONDO'
Set-TextValue 45 3 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 45 4 '1.21'
Set-TextValue 45 5 '  +7.65%  '
Set-TextValue 46 4 '0.761'
Set-TextValue 46 5 '  +1.61%  '
Set-TextValue 47 5 '  +6.21%  '
Set-TextValue 48 2 'EnergySwap'
Set-TextValue 48 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 48 4 '23.40'
Set-TextValue 48 5 '  +4.16%  '
Set-TextValue 49 2 'Cosmos'
Set-TextValue 49 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 49 4 '6.99'
Set-TextValue 49 5 '  +5.98%  '
Set-TextValue 50 4 '23.20'
Set-TextValue 50 5 '  +11.73%  '
Set-TextValue 51 4 '2.45'
Set-TextValue 51 5 '  +14.10%  '
